# Add two new DQL test cases (dqlc3_007, dqlc3_008) covering getByKey on a
# non-primary-key column, per commit message:
#   "add case about query non-key but operation is get and cause issu"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: dqlc3_007 - getByKey output for non-key int column query ---
$ws.Range("A8").Value = "dqlc3_007"
$ws.Range("C8").Value = "getByKey输出 - 非主键int型查询"
$ws.Range("D8").Value = "Explain-key"
$ws.Range("F8").Value = "schema17"
$ws.Range("G8").Value = "qc3_value1"
$ws.Range("I8").Value = "select age from `$schema17 where age=18"
$ws.Range("J8").Value = "src/test/resources/io.dingodb.test/testdata/cases/dql/casegroup3/expectedresult/queryc3_007.csv"
$ws.Range("K8").Value = "csv_containsAll"

# --- Row 9: dqlc3_008 - getByKey output for non-key varchar column query ---
$ws.Range("A9").Value = "dqlc3_008"
$ws.Range("C9").Value = "getByKey输出 - 非主键varchar型查询"
$ws.Range("D9").Value = "Explain-key"
$ws.Range("F9").Value = "schema17"
$ws.Range("G9").Value = "qc3_value1"
$ws.Range("I9").Value = "select name from `$schema17 where name='zhangsan'"
$ws.Range("J9").Value = "src/test/resources/io.dingodb.test/testdata/cases/dql/casegroup3/expectedresult/queryc3_008.csv"
$ws.Range("K9").Value = "csv_containsAll"

# Column B ("Testable") uses a distinct bold-ish style (same as the rest of
# the data rows) - copy that formatting from an existing row, then set the
# value, so the new cells match the existing "y/n" list-validated column.
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B8").Value = "y"
$ws.Range("B9").Value = "y"

# Match the saved view state: scrolled back to column A, selection on J9.
$ws.Range("J9").Select()
